# Trade #5 closed at 2026-02-17 19:56:07 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade #5 (MarketMaking, UP, entry 0.58, exit 0.55).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.75   # Current Capital
$summary.Range("B4").Value = -0.25     # Total P&L $
$summary.Range("B5").Value = -1        # Total P&L %
$summary.Range("B6").Value = 5         # Total Trades
$summary.Range("B8").Value = 4         # Losing Trades
$summary.Range("B9").Value = 20        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.75   # Capital
$status.Range("D5").Value = 5       # Trades
$status.Range("E5").Value = -0.25   # P&L $
$status.Range("F5").Value = -0.25   # P&L %
$status.Range("G5").Value = 20      # Win Rate %

# ---------------------------------------------------------------------------
# All Trades + MarketMaking sheets - append trade #5
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(6, 1).Value = 5

    # The Date column holds plain text like "2026-02-17" in every other row
    # (not an actual date value). Force text so COM doesn't auto-convert the
    # string into a date serial number, then drop the explicit number-format
    # style so the cell matches the unformatted (default-style) neighbours.
    $dateCell = $ws.Cells.Item(6, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $ws.Cells.Item(6, 3).Value = "19:56:00"
    $ws.Cells.Item(6, 4).Value = "MarketMaking"
    $ws.Cells.Item(6, 5).Value = "UP"
    $ws.Cells.Item(6, 6).Value = 0.58
    $ws.Cells.Item(6, 7).Value = 0.55
    $ws.Cells.Item(6, 8).Value = "CLOSED"
    $ws.Cells.Item(6, 9).Value = -5.1724
    $ws.Cells.Item(6, 10).Value = -0.03
    $ws.Cells.Item(6, 11).Value = 99.75
    $ws.Cells.Item(6, 12).Value = 0
    $ws.Cells.Item(6, 13).Value = 0
    $ws.Cells.Item(6, 14).Value = 0.6
    $ws.Cells.Item(6, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(6, 16).Value = "early_exit"
    $ws.Cells.Item(6, 17).Value = 0.13
}
